$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 42.05115733333333
$ws.Range("H2").Value = 126.153472
$ws.Range("I2").Value = 0.1594435451835853
$ws.Range("J2").Value = 0.1594435451835853
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.070922
$ws.Range("N2").Value = 0.212766
$ws.Range("O2").Value = 0.02446124354254487
$ws.Range("P2").Value = 0.02446124354254487
$ws.Range("Q2").Value = 2.982352180394666
$ws.Range("R2").Value = 26.841169623552
$ws.Range("S2").Value = 0.003900187390022437
$ws.Range("T2").Value = 0.003900187390022438
# Row 3
$ws.Range("G3").Value = 42.05115733333333
$ws.Range("H3").Value = 126.153472
$ws.Range("I3").Value = 0.1594435451835853
$ws.Range("J3").Value = 0.1594435451835853
$ws.Range("M3").Value = 0.1465633333333334
$ws.Range("O3").Value = 0.05055020150410101
$ws.Range("P3").Value = 0.05055020150410101
$ws.Range("Q3").Value = 6.163157789297778
$ws.Range("R3").Value = 55.46842010368
$ws.Range("S3").Value = 0.008059903337558471
$ws.Range("T3").Value = 0.008059903337558471
# Row 4
$ws.Range("G4").Value = 42.05115733333333
$ws.Range("H4").Value = 126.153472
$ws.Range("I4").Value = 0.1594435451835853
$ws.Range("J4").Value = 0.1594435451835853
$ws.Range("M4").Value = 2.681876666666667
$ws.Range("N4").Value = 8.045630000000001
$ws.Range("O4").Value = 0.9249885549533541
$ws.Range("P4").Value = 0.9249885549533541
$ws.Range("Q4").Value = 112.7760176585956
$ws.Range("R4").Value = 1014.98415892736
$ws.Range("S4").Value = 0.1474834544560044
$ws.Range("T4").Value = 0.1474834544560044
# Row 5
$ws.Range("G5").Value = 57.66057933333332
$ws.Range("I5").Value = 0.2186291119973147
$ws.Range("J5").Value = 0.2186291119973148
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.070922
$ws.Range("N5").Value = 0.212766
$ws.Range("O5").Value = 0.02446124354254487
$ws.Range("P5").Value = 0.02446124354254487
$ws.Range("Q5").Value = 4.089403607478665
$ws.Range("R5").Value = 36.804632467308
$ws.Range("S5").Value = 0.005347939954056635
$ws.Range("T5").Value = 0.005347939954056636
# Row 6
$ws.Range("G6").Value = 57.66057933333332
$ws.Range("I6").Value = 0.2186291119973147
$ws.Range("J6").Value = 0.2186291119973148
$ws.Range("M6").Value = 0.1465633333333334
$ws.Range("O6").Value = 0.05055020150410101
$ws.Range("P6").Value = 0.05055020150410101
$ws.Range("Q6").Value = 8.450926709024444
$ws.Range("S6").Value = 0.01105174566612693
$ws.Range("T6").Value = 0.01105174566612693
# Row 7
$ws.Range("G7").Value = 57.66057933333332
$ws.Range("I7").Value = 0.2186291119973147
$ws.Range("J7").Value = 0.2186291119973148
$ws.Range("M7").Value = 2.681876666666667
$ws.Range("N7").Value = 8.045630000000001
$ws.Range("O7").Value = 0.9249885549533541
$ws.Range("P7").Value = 0.9249885549533541
$ws.Range("Q7").Value = 154.6385623005489
$ws.Range("R7").Value = 1391.74706070494
$ws.Range("S7").Value = 0.2022294263771312
$ws.Range("T7").Value = 0.2022294263771312
# Row 8
$ws.Range("G8").Value = 99.15200299999999
$ws.Range("H8").Value = 297.456009
$ws.Range("I8").Value = 0.3759503393701321
$ws.Range("J8").Value = 0.3759503393701321
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.070922
$ws.Range("N8").Value = 0.212766
$ws.Range("O8").Value = 0.02446124354254487
$ws.Range("P8").Value = 0.02446124354254487
$ws.Range("Q8").Value = 7.032058356765999
$ws.Range("R8").Value = 63.288525210894
$ws.Range("S8").Value = 0.009196212811235198
$ws.Range("T8").Value = 0.009196212811235198
# Row 9
$ws.Range("G9").Value = 99.15200299999999
$ws.Range("H9").Value = 297.456009
$ws.Range("I9").Value = 0.3759503393701321
$ws.Range("J9").Value = 0.3759503393701321
$ws.Range("M9").Value = 0.1465633333333334
$ws.Range("O9").Value = 0.05055020150410101
$ws.Range("P9").Value = 0.05055020150410101
$ws.Range("Q9").Value = 14.53204806635667
$ws.Range("R9").Value = 130.78843259721
$ws.Range("S9").Value = 0.01900436541069534
$ws.Range("T9").Value = 0.01900436541069534
# Row 10
$ws.Range("G10").Value = 99.15200299999999
$ws.Range("H10").Value = 297.456009
$ws.Range("I10").Value = 0.3759503393701321
$ws.Range("J10").Value = 0.3759503393701321
$ws.Range("M10").Value = 2.681876666666667
$ws.Range("N10").Value = 8.045630000000001
$ws.Range("O10").Value = 0.9249885549533541
$ws.Range("P10").Value = 0.9249885549533541
$ws.Range("Q10").Value = 265.9134432989634
$ws.Range("R10").Value = 2393.22098969067
$ws.Range("S10").Value = 0.3477497611482015
$ws.Range("T10").Value = 0.3477497611482016
# Row 11
$ws.Range("G11").Value = 64.87322933333333
$ws.Range("H11").Value = 194.619688
$ws.Range("I11").Value = 0.2459770034489679
$ws.Range("J11").Value = 0.2459770034489679
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.070922
$ws.Range("N11").Value = 0.212766
$ws.Range("O11").Value = 0.02446124354254487
$ws.Range("P11").Value = 0.02446124354254487
$ws.Range("Q11").Value = 4.600939170778666
$ws.Range("R11").Value = 41.408452537008
$ws.Range("S11").Value = 0.006016903387230603
$ws.Range("T11").Value = 0.006016903387230604
# Row 12
$ws.Range("G12").Value = 64.87322933333333
$ws.Range("H12").Value = 194.619688
$ws.Range("I12").Value = 0.2459770034489679
$ws.Range("J12").Value = 0.2459770034489679
$ws.Range("M12").Value = 0.1465633333333334
$ws.Range("O12").Value = 0.05055020150410101
$ws.Range("P12").Value = 0.05055020150410101
$ws.Range("Q12").Value = 9.508036735191112
$ws.Range("R12").Value = 85.57233061672001
$ws.Range("S12").Value = 0.01243418708972028
$ws.Range("T12").Value = 0.01243418708972028
# Row 13
$ws.Range("G13").Value = 64.87322933333333
$ws.Range("H13").Value = 194.619688
$ws.Range("I13").Value = 0.2459770034489679
$ws.Range("J13").Value = 0.2459770034489679
$ws.Range("M13").Value = 2.681876666666667
$ws.Range("N13").Value = 8.045630000000001
$ws.Range("O13").Value = 0.9249885549533541
$ws.Range("P13").Value = 0.9249885549533541
$ws.Range("Q13").Value = 173.9820000403822
$ws.Range("R13").Value = 1565.83800036344
$ws.Range("S13").Value = 0.227525912972017
$ws.Range("T13").Value = 0.227525912972017
